# Revert "Added some classes to C# server, and modified documents"
#
# Slide 1, shape 6:  "Utilities" / "(Include Win 32 API)"  (2 paragraphs)
#                     -> single paragraph "Win 32 API"
# Slide 3, shape 11: "Page Setter"   -> runs "Page " + "Setter"
# Slide 3, shape 12: "Screen Getter" -> runs "Screen " + "Getter"
# Slide 5, shape 7:  "Utilities"     -> "Win 32 API" (simple run text replace)
# Slide 5, shape 11: "Page Setter"   -> runs "Page " + "Setter"
# Slide 5, shape 12: "Screen Getter" -> runs "Screen " + "Getter"

$p = $ppt.ActivePresentation

function Merge-IntoWin32Api($shape) {
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf([char]13)
    if ($idx -ge 0) {
        # Remove the first paragraph (including its paragraph mark) so the
        # remaining paragraph keeps its own pPr/endParaRPr intact.
        $firstPara = $tr.Characters(1, $idx + 1)
        $firstPara.Delete()
    }
    $tr.Text = "Win 32 API"
}

function Split-Run($shape, $firstPart, $secondPart) {
    $tr = $shape.TextFrame.TextRange
    $tailLen = $secondPart.Length
    $total = $tr.Text.Length
    $tail = $tr.Characters($total - $tailLen + 1, $tailLen)
    # Re-assigning the same text on the sub-range forces PowerPoint to
    # split the run in two while copying the original run's properties.
    $tail.Text = $secondPart
}

# --- Slide 1 ---
$s1 = $p.Slides.Item(1)
Merge-IntoWin32Api $s1.Shapes.Item(6)

# --- Slide 3 ---
$s3 = $p.Slides.Item(3)
Split-Run $s3.Shapes.Item(11) "Page " "Setter"
Split-Run $s3.Shapes.Item(12) "Screen " "Getter"

# --- Slide 5 ---
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(7).TextFrame.TextRange.Text = "Win 32 API"
Split-Run $s5.Shapes.Item(11) "Page " "Setter"
Split-Run $s5.Shapes.Item(12) "Screen " "Getter"
